$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM score values for rows 2-21 (residues), columns B-K (positions 0-9).
# Each line below is one sheet row of 10 comma-separated numbers (cols B..K).
$rowsText = @(
    "-16.85270018050934,0.9130844653505339,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,2.398950962453236,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,0.9219903771057719,-16.85270018050934,-16.85270018050934,1.716120640901276,-16.85270018050934,1.716194843875742,-16.85270018050934,2.447449021711271,-16.85270018050934",
    "-16.85270018050934,2.008571534711804,-16.85270018050934,-16.85270018050934,-16.85270018050934,2.101885224211973,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "3.033825496929874,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,3.000031998584574,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "3.561944800968407,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,1.551131550268683,-16.85270018050934,2.215534906685878",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,2.539434830852328,-16.85270018050934,2.465571198068851,-16.85270018050934,-16.85270018050934,-16.85270018050934,1.337187101200288",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,0.9522941922142805,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,2.266790656817293,1.61054079208639",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,2.051666897832817",
    "-16.85270018050934,-16.85270018050934,4.321916514272325,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,2.320247324256041,-16.85270018050934",
    "-16.85270018050934,1.845899367297521,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,0.4386234300503786,0.7641429729010896,1.161303939052872,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,0.5246883029301406,0.8227951889129053,1.32312709961201,-16.85270018050934",
    "-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,-16.85270018050934,1.865635354987988,2.114551412643439,-16.85270018050934,-16.85270018050934",
    "-16.85270018050934,2.4627877391258,-16.85270018050934,-16.85270018050934,4.063013652500737,-16.85270018050934,2.248436840227046,1.999088504224348,-16.85270018050934,2.492941123950489",
    "-16.85270018050934,1.627920432115021,-16.85270018050934,2.087845702122877,-16.85270018050934,3.348190045428195,2.465791877661494,-16.85270018050934,-16.85270018050934,-16.85270018050934"
)

$numRows = $rowsText.Count
$numCols = 10
$arr = New-Object "object[,]" $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    $parts = $rowsText[$r].Split(",")
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r,$c] = [double]$parts[$c]
    }
}

$ws.Range("B2:K21").Value = $arr
